$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D26").Value = "2021 인공지능 경진대회 참가기"
$ws.Range("E26").Value = "https://blog.est.ai/2021/12/2021aicontest/"

$ws.Range("D28").Value = "Let's do MuJoCo - 1. Mujoco, mujoco-py 설치"
$ws.Range("E28").Value = "https://ropiens.tistory.com/169"

$wb.Save()
